# Update G418-identified strains: mark marker_1 (column J) as "G418"
# for the TDY1210 / CNAG_00440 rows (17-19), and move the active
# selection to J3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J17:J19").Value = "G418"

$ws.Range("J3").Select()
